# Fruta / hortaliza, semanal
# Inserts the newest week's data (Primera + Segunda grades) for
# "Feria Lagunitas de Puerto Montt - Pomelo" right above the existing
# historical rows (which start at row 273), pushing the remaining
# historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 273 (formats inherited from the row above),
# shifting old rows 273..402 down to 275..404.
$ws.Rows.Item(273).Resize(2).Insert()

# New row 273: Primera grade, newest week
$ws.Cells.Item(273,1).Value2  = 4
$ws.Cells.Item(273,2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(273,3).Value2  = "Los Lagos"
$ws.Cells.Item(273,4).Value2  = 44839
$ws.Cells.Item(273,5).Value2  = 10
$ws.Cells.Item(273,6).Value2  = "Fruta"
$ws.Cells.Item(273,7).Value2  = 100102
$ws.Cells.Item(273,8).Value2  = "Cítricos"
$ws.Cells.Item(273,9).Value2  = 100102006
$ws.Cells.Item(273,10).Value2 = "Pomelo"
$ws.Cells.Item(273,11).Value2 = "Start Ruby"
$ws.Cells.Item(273,12).Value2 = "Primera"
$ws.Cells.Item(273,13).Value2 = 20
$ws.Cells.Item(273,14).Value2 = 13000
$ws.Cells.Item(273,15).Value2 = 14000
$ws.Cells.Item(273,16).Value2 = 13500
$ws.Cells.Item(273,17).Value2 = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(273,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(273,19).Value2 = 964
$ws.Cells.Item(273,20).Value2 = 14

# New row 274: Segunda grade, newest week
$ws.Cells.Item(274,1).Value2  = 4
$ws.Cells.Item(274,2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(274,3).Value2  = "Los Lagos"
$ws.Cells.Item(274,4).Value2  = 44839
$ws.Cells.Item(274,5).Value2  = 10
$ws.Cells.Item(274,6).Value2  = "Fruta"
$ws.Cells.Item(274,7).Value2  = 100102
$ws.Cells.Item(274,8).Value2  = "Cítricos"
$ws.Cells.Item(274,9).Value2  = 100102006
$ws.Cells.Item(274,10).Value2 = "Pomelo"
$ws.Cells.Item(274,11).Value2 = "Start Ruby"
$ws.Cells.Item(274,12).Value2 = "Segunda"
$ws.Cells.Item(274,13).Value2 = 10
$ws.Cells.Item(274,14).Value2 = 12000
$ws.Cells.Item(274,15).Value2 = 12000
$ws.Cells.Item(274,16).Value2 = 12000
$ws.Cells.Item(274,17).Value2 = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(274,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(274,19).Value2 = 857
$ws.Cells.Item(274,20).Value2 = 14
